# Add Claas.Rostock@DNV.com to the author/contact line(s) that currently
# read "Siegfried.Eisinger@DNV.com & Jorge.Luis.Mendez@DNV.com".
#
# For every such occurrence in the deck we:
#   1. Turn " & " into ", " (so the line reads "A, B").
#   2. Re-apply the hyperlink on Jorge's e-mail (it keeps its own run).
#   3. Append " and Claas.Rostock@DNV.com" and hyperlink the new e-mail.
#
# The very first occurrence (slide 1) gets the e-mail as a single run.
# The later occurrence (slide 10) was re-typed by hand afterwards and ended
# up split into three runs ("Claas.Rostock@" / "DNV" / ".com") -- we
# reproduce that same run layout there.

$p = $ppt.ActivePresentation

$oldSep    = " & "
$jorgeName = "Jorge.Luis.Mendez@DNV.com"
$oldSpan   = $oldSep + $jorgeName

$claasName = "Claas.Rostock@DNV.com"
$claasMail = "mailto:" + $claasName
$jorgeMail = "mailto:" + $jorgeName

$occurrence = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)

        if (-not $shp.HasTextFrame) { continue }
        if (-not $shp.TextFrame.HasText) { continue }

        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text

        if ($full.IndexOf($oldSpan) -lt 0) { continue }

        $occurrence = $occurrence + 1

        # --- Step 1: replace " & Jorge.Luis.Mendez@DNV.com" with a clean,
        # non-hyperlinked ", Jorge.Luis.Mendez@DNV.com and Claas.Rostock@DNV.com"
        # (selecting from inside the plain " & " run keeps the whole
        # replacement free of any inherited hyperlink formatting).
        $spanStart = $full.IndexOf($oldSpan) + 1
        $span = $tr.Characters($spanStart, $oldSpan.Length)
        $newSpan = ", " + $jorgeName + " and " + $claasName
        $span.Text = $newSpan

        $full2 = $tr.Text

        # --- Step 2: re-hyperlink Jorge's e-mail.
        $jorgeStart = $full2.IndexOf($jorgeName) + 1
        $jorgeRange = $tr.Characters($jorgeStart, $jorgeName.Length)
        $jorgeRange.ActionSettings(1).Hyperlink.Address = $jorgeMail

        # --- Step 3: hyperlink Claas's e-mail.
        $claasStart = $full2.IndexOf($claasName) + 1

        if ($occurrence -eq 1) {
            # First occurrence: one single run for the whole address.
            $claasRange = $tr.Characters($claasStart, $claasName.Length)
            $claasRange.ActionSettings(1).Hyperlink.Address = $claasMail
        } else {
            # Later occurrence: split into three runs, as it ended up being
            # typed/corrected by hand in the source deck.
            $part1 = "Claas.Rostock@"
            $part2 = "DNV"
            $part3 = ".com"

            $r1 = $tr.Characters($claasStart, $part1.Length)
            $r1.ActionSettings(1).Hyperlink.Address = $claasMail

            $r2start = $claasStart + $part1.Length
            $r2 = $tr.Characters($r2start, $part2.Length)
            $r2.ActionSettings(1).Hyperlink.Address = $claasMail

            $r3start = $r2start + $part2.Length
            $r3 = $tr.Characters($r3start, $part3.Length)
            $r3.ActionSettings(1).Hyperlink.Address = $claasMail
        }
    }
}

Write-Host ("Updated " + $occurrence + " occurrence(s).")
